# Generate Report for Handback
# Update "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 73e28401-ac3f-4dda-8550-b29fa5410a52 row across all three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest HO Xliff Generate Date (column G) for row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-16 02:43:42"

# --- zh-cn sheet: Correspond Handoff/Handback Datetime for row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-16 02:43:37"
$wsZhCn.Range("K3").Value = "2016-08-16 02:43:54"

# --- de-de sheet: Correspond Handoff/Handback Datetime for row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-16 02:43:42"
$wsDeDe.Range("K3").Value = "2016-08-16 02:44:03"
